$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two duplicate trailing rows (12 and 13) first so the row
# numbers used below line up with the final 10-row layout.
$ws.Rows("12:13").Delete()

# --- Column D (Filename) : swap in the 10 new TRY-file names ---
$ws.Range("D2").Value = "TRY2015_515220074856_Jahr.dat"
$ws.Range("D3").Value = "TRY2015_509319069572_Jahr.dat"
$ws.Range("D4").Value = "TRY2015_507755060854_Jahr.dat"
$ws.Range("D5").Value = "TRY2015_513148094876_Jahr.dat"
$ws.Range("D6").Value = "TRY2015_525153133939_Jahr.dat"
$ws.Range("D7").Value = "TRY2015_510342136998_Jahr.dat"
$ws.Range("D8").Value = "TRY2015_480091078440_Jahr.dat"
$ws.Range("D9").Value = "TRY2015_535485100234_Jahr.dat"
$ws.Range("D10").Value = "TRY2015_481593115227_Jahr.dat"
$ws.Range("D11").Value = "TRY2015_522733105384_Jahr.dat"

# --- Column A (PLZ) : replace with PLZ codes of 10 example cities ---
# Enter as real numbers first ...
$ws.Range("A2").Value = 44137
$ws.Range("A3").Value = 50667
$ws.Range("A4").Value = 52064
$ws.Range("A5").Value = 34117
$ws.Range("A6").Value = 10115
$ws.Range("A8").Value = 80331
$ws.Range("A9").Value = 20354
$ws.Range("A10").Value = 79100
$ws.Range("A11").Value = 38100

# ... then switch the whole column to Text format (keeps the numbers as
# numbers but displayed/stored the way the workbook was saved) ...
$ws.Range("A2:A11").NumberFormat = "@"

# ... Dresden's PLZ keeps its leading zero, so it has to be entered as text.
$ws.Range("A7").Value = "01067"

# --- cosmetic touch-ups to match the saved worksheet state ---
$ws.Columns("D").ColumnWidth = 38.19

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A2").Select()
